$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "55.811.50"
$ws.Range("E2").Value = "  -1.01%  "
$ws.Range("D3").Value = "2.375.25"
$ws.Range("E3").Value = "  -4.84%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'477.37"
$ws.Range("E5").Value = "  -2.35%  "
$ws.Range("D6").Value = "'146.99"
$ws.Range("E6").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.32%  "
$ws.Range("E8").Value = "  -2.26%  "
$ws.Range("D9").Value = "2.374.92"
$ws.Range("E9").Value = "  -5.66%  "
$ws.Range("D10").Value = "'0.0973"
$ws.Range("E10").Value = "  -0.38%  "
$ws.Range("D11").Value = "'5.42"
$ws.Range("E11").Value = "  -6.26%  "
$ws.Range("D12").Value = "'0.323"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "2.789.71"
$ws.Range("E14").Value = "  -5.02%  "
$ws.Range("D15").Value = "55.905.14"
$ws.Range("E15").Value = "  -0.90%  "
$ws.Range("D16").Value = "'20.33"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("E17").Value = "  -3.48%  "
$ws.Range("D18").Value = "2.378.49"
$ws.Range("E18").Value = "  -5.32%  "
$ws.Range("D19").Value = "'4.59"
$ws.Range("E19").Value = "  +1.38%  "
$ws.Range("D20").Value = "'315.40"
$ws.Range("E20").Value = "  -1.73%  "
$ws.Range("D21").Value = "'9.70"
$ws.Range("E21").Value = "  -5.21%  "
$ws.Range("D22").Value = "'0.997"
$ws.Range("E22").Value = "  -0.08%  "
$ws.Range("D23").Value = "'5.66"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("D24").Value = "'56.72"
$ws.Range("E24").Value = "  -3.39%  "
$ws.Range("E25").Value = "  +0.32%  "
$ws.Range("D26").Value = "'0.395"
$ws.Range("E26").Value = "  -4.15%  "
$ws.Range("E27").Value = "  -5.78%  "
$ws.Range("D28").Value = "2.484.83"
$ws.Range("E28").Value = "  -5.12%  "
$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  -5.12%  "
$ws.Range("D30").Value = "0.0₃0768"
$ws.Range("E30").Value = "  -3.20%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("D32").Value = "'147.04"
$ws.Range("E32").Value = "  -1.11%  "
$ws.Range("D33").Value = "'18.06"
$ws.Range("E33").Value = "  -1.24%  "
$ws.Range("E34").Value = "  -1.80%  "
$ws.Range("D35").Value = "'5.06"
$ws.Range("E35").Value = "  -3.13%  "
$ws.Range("E36").Value = "  -3.94%  "
$ws.Range("D37").Value = "'3.58"
$ws.Range("E37").Value = "  -4.29%  "
$ws.Range("D38").Value = "'0.836"
$ws.Range("E38").Value = "  -3.88%  "
$ws.Range("D39").Value = "'33.41"
$ws.Range("E39").Value = "  -2.41%  "
$ws.Range("E40").Value = "  +0.39%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -4.70%  "
$ws.Range("D43").Value = "'0.0536"
$ws.Range("E43").Value = "  -4.09%  "
$ws.Range("D44").Value = "'0.0947"
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "'0.581"
$ws.Range("E45").Value = "  -6.16%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'254.36"
$ws.Range("E47").Value = "  -2.56%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "'0.0223"
$ws.Range("E48").Value = "  -2.60%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.53"
$ws.Range("E49").Value = "  -6.94%  "
$ws.Range("D50").Value = "'16.90"
$ws.Range("E50").Value = "  -4.51%  "
$ws.Range("D51").Value = "1.779.71"
$ws.Range("E51").Value = "  -7.15%  "
